# Apply odds updates to Sheet1 of the FlashScore "Jogos da Semana" workbook.
# Each assignment below mirrors a single <v> value change from the source
# OOXML diff, addressed by its worksheet cell reference.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2.5
$ws.Range("I2").Value = 3.1
$ws.Range("L2").Value = 4
$ws.Range("Z2").Value = 23
$ws.Range("AA2").Value = 23
$ws.Range("BB2").Value = 351

# Row 3
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5

# Row 4
$ws.Range("G4").Value = 3.75
$ws.Range("H4").Value = 2.88
$ws.Range("I4").Value = 2.2
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("O4").Value = 1.53
$ws.Range("P4").Value = 2.38
$ws.Range("AA4").Value = 41
$ws.Range("AI4").Value = 9
$ws.Range("AN4").Value = 5.5
$ws.Range("AP4").Value = 41

# Row 5
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 3
$ws.Range("Q5").Value = 2.15
$ws.Range("R5").Value = 1.67

# Row 6
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 10

# Row 8
$ws.Range("G8").Value = 1.42
$ws.Range("H8").Value = 4.2
$ws.Range("I8").Value = 8
$ws.Range("J8").Value = 1.95
$ws.Range("L8").Value = 7
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 11
$ws.Range("U8").Value = 2
$ws.Range("V8").Value = 1.73
$ws.Range("Y8").Value = 8.5
$ws.Range("Z8").Value = 9.5
$ws.Range("AD8").Value = 8
$ws.Range("AE8").Value = 19
$ws.Range("AK8").Value = 81
$ws.Range("AP8").Value = 19
$ws.Range("BA8").Value = 151

# Row 9
$ws.Range("G9").Value = 1.42
$ws.Range("H9").Value = 4.1
$ws.Range("I9").Value = 8.5

# Row 11
$ws.Range("G11").Value = 1.5
$ws.Range("H11").Value = 3.8
$ws.Range("I11").Value = 7.5
$ws.Range("J11").Value = 2.1
$ws.Range("K11").Value = 2.2
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("U11").Value = 2.2
$ws.Range("V11").Value = 1.62
$ws.Range("AJ11").Value = 23
$ws.Range("AW11").Value = 8

# Row 14
$ws.Range("G14").Value = 6.5
$ws.Range("H14").Value = 4.33
$ws.Range("I14").Value = 1.48
$ws.Range("J14").Value = 7
$ws.Range("L14").Value = 2.05
$ws.Range("M14").Value = 1.06
$ws.Range("N14").Value = 10
$ws.Range("O14").Value = 1.29
$ws.Range("P14").Value = 3.5
$ws.Range("Q14").Value = 1.98
$ws.Range("R14").Value = 1.88
$ws.Range("W14").Value = 15
$ws.Range("X14").Value = 34
$ws.Range("Y14").Value = 21
$ws.Range("Z14").Value = 81
$ws.Range("AD14").Value = 8.5
$ws.Range("AK14").Value = 9.5
$ws.Range("AN14").Value = 8
$ws.Range("AO14").Value = 41
$ws.Range("AQ14").Value = 151
$ws.Range("AU14").Value = 9.5
$ws.Range("AW14").Value = 3.25
